# ErroresDelSistema.xlsx - "Carga del administrador de Catalogos" commit
#
# Adds 7 new error/info rows (35-41, 1-based data rows -> sheet rows 36-42)
# describing the new "Catalogos" (catalog admin) feature to the Errores
# sheet, and marks the formatting cell E43 with the existing underline
# style used elsewhere on the sheet. The Script sheet's formulas
# automatically recompute the generated INSERT statements from these new
# values, so nothing needs to be touched there directly.

$wb = $excel.ActiveWorkbook
$wsErrores = $wb.Worksheets.Item("Errores")
$wsScript  = $wb.Worksheets.Item("Script")

# --- Row 36 (id 35): catalog update success ------------------------------
$wsErrores.Range("C36").Value = "Se a actualizado correctamente el catálogo"
$wsErrores.Range("D36").Value = "INFO"
$wsErrores.Range("E36").Value = "Éxito al modificar el catálogo"

# --- Row 37 (id 36): session restart warning ------------------------------
$wsErrores.Range("C37").Value = "Para que el cambio tenga efecto deberá reiniciar su sesión."
$wsErrores.Range("D37").Value = "WARN"
$wsErrores.Range("E37").Value = "Éxito al modificar el catálogo"

# --- Row 38 (id 37): catalog update error ---------------------------------
$wsErrores.Range("C38").Value = "Se ha presentado un problema al modificar el catálogo. Intente nuevamente. Si el problema persiste contacte con su administrador"
$wsErrores.Range("D38").Value = "ERROR"
$wsErrores.Range("E38").Value = "Error al modificar el catálogo"

# --- Row 39 (id 38): catalog create success -------------------------------
$wsErrores.Range("C39").Value = "El catálogo `$ ha sido guardado correctamente, junto con el item `$."
$wsErrores.Range("D39").Value = "INFO"
$wsErrores.Range("E39").Value = "Éxito al crear el catálogo"

# --- Row 40 (id 39): catalog create error ---------------------------------
$wsErrores.Range("C40").Value = "El catálogo `$ no ha sido guardado correctamente, junto con el item `$. Si el problema perisite contacte con su administrador"
$wsErrores.Range("D40").Value = "ERROR"
$wsErrores.Range("E40").Value = "Error al crear el catálogo"

# --- Row 41 (id 40): catalog create fatal error ---------------------------
$wsErrores.Range("C41").Value = "Ha ocurrido un error inesperado al guardar el catálogo. Si el problema persiste contacte con su administrador"
$wsErrores.Range("D41").Value = "FATAL"
$wsErrores.Range("E41").Value = "Error al crear el catálogo"

# --- Row 42 (id 41): view retrieval fatal error ---------------------------
$wsErrores.Range("C42").Value = "Se ha presentado un error al tratar de recuperar los catálogos actuales. Si el problema persiste contacte con su administrador"
$wsErrores.Range("D42").Value = "FATAL"
$wsErrores.Range("E42").Value = "Error al recuperar la vista"

# --- Row 43 (id 42): stray formatted-but-empty cell, matches the already
#     present underline style (s="5", same as used on A37/F25/G25) -------
$wsErrores.Range("E43").Font.Underline = $true

# --- Selection / active-sheet state ---------------------------------------
# Errores keeps a selection over the newly added rows (not the active tab
# any more)...
$wsErrores.Range("C36:C42").Select()

# ...while Script becomes the active tab, selected over the corresponding
# generated-formula rows.
$wsScript.Activate()
$wsScript.Range("B36:B42").Select()
